$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are stored as text (the sheet has no
# number formatting and percentages carry literal padding spaces). Some new
# Price values (e.g. "1.00", "211.41") look numeric, so a bare .Value
# assignment would make Excel auto-convert them to numbers and silently
# drop formatting (e.g. "154.80" -> 154.8). Prefixing those with a leading
# single-quote forces Excel to keep them as literal text, matching the
# original cell type.

$ws.Range("D2").Value = '27.901.94'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.632.77'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''211.41'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''23.41'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").Value = '''0.257'
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '''0.0882'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '1.865.27'
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = '1.639.63'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").Value = '''65.46'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '27.913.41'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '''229.16'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = '''4.33'
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '''154.80'
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("D26").Value = '''6.88'
$ws.Range("E26").Value = '  -1.01%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = '''15.53'
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").Value = '''3.41'
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").Value = '''3.12'
$ws.Range("E33").Value = '  +1.47%  '
$ws.Range("D34").Value = '1.393.26'
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("E36").Value = '  +10.02%  '
$ws.Range("E37").Value = '  -1.09%  '
$ws.Range("E38").Value = '  +1.21%  '
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").Value = '''0.849'
$ws.Range("E40").Value = '  -2.97%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''65.74'
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''1.82'
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '''5.42'
$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("D46").Value = '1.775.06'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("D48").Value = '''88.74'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = '''7.65'
$ws.Range("E51").Value = '  +1.10%  '
